$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'55.493.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.46%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.952.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.98%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'487.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.65%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'131.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'2.949.85"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.99%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.418"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.35%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'7.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.45%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.100"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -8.45%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "  -9.93%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'3.446.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.48%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'24.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.98%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'55.344.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.74%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.942.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.50%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "  -8.26%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "  -2.67%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'12.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.22%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'7.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.13%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'313.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.83%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.460"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -9.37%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'59.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -13.49%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "  -5.28%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.0₃0841"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -10.82%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'6.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.32%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'6.53"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.47%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.49%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "  -10.18%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'19.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -10.72%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'147.95"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.55%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'4.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -10.08%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'5.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -9.23%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "  -8.51%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'23.62"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.07%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.0644"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.81%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.974.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.18%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'35.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -11.62%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.39%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.629"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.21%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.78%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'3.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.91%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.124.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.58%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0231"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.87%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'18.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.49%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'5.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.50%  "
$ws.Range("E51").Style = "Normal"
